$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 9 timestamp (float precision correction) ---
$ws.Cells.Item(9, 1).Value = 45730.4355013426

# --- New row 10 ---
$ws.Cells.Item(10, 1).Value = 45730.50495601852
$ws.Cells.Item(10, 2).Value = 'fff'
$ws.Cells.Item(10, 3).Value = 25
$ws.Cells.Item(10, 4).Value = 'Male'
$ws.Cells.Item(10, 5).Value = 50
$ws.Cells.Item(10, 6).Value = 1.5
$ws.Cells.Item(10, 7).Value = 22.22
$ws.Cells.Item(10, 8).Value = 'Normal weight - Maintain a balanced diet and exercise.'
$ws.Cells.Item(10, 9).Value = 'Veg'
$ws.Cells.Item(10, 10).Value = 'vitamin_C'
$ws.Cells.Item(10, 11).Value = @'

Recommendations for vitamin_C Deficiency:
Acerola, (west indian cherry), raw
Acerola juice, raw
Guavas, common, raw
Peppers, hot chili, green, raw
Peppers, sweet, yellow, raw
Mango, Ataulfo, peeled, raw
Currants, european black, raw
Kiwifruit, ZESPRI SunGold, raw
Peppers, bell, orange, raw
Drumstick pods, raw
Pokeberry shoots, (poke), raw
Lemon peel, raw
Peppers, bell, yellow, raw
Peppers, bell, red, raw
Orange peel, raw
Longans, raw
Peppers, bell, green, raw
Balsam-pear (bitter gourd), pods, raw
Peppers, sweet, red, raw
Litchis, raw
Broccoli, flower clusters, raw
Broccoli, leaves, raw
Broccoli, stalks, raw
Jujube, raw
Kiwifruit, green, raw
Peppers, sweet, green, raw
Persimmons, native, raw
Taro, tahitian, raw
Pummelo, raw
Strawberries, raw
Oranges, raw, with peel
Papayas, raw
Strawberries, raw
Broccoli, raw
Kiwifruit (kiwi), green, peeled, raw
Mustard spinach, (tendergreen), raw
Pineapple, raw
Cauliflower, green, raw
'@
$ws.Cells.Item(10, 1).NumberFormat = $ws.Cells.Item(9, 1).NumberFormat

# --- New row 11 ---
$ws.Cells.Item(11, 1).Value = 45732.90905953704
$ws.Cells.Item(11, 2).Value = 'eddc'
$ws.Cells.Item(11, 3).Value = 25
$ws.Cells.Item(11, 4).Value = 'Male'
$ws.Cells.Item(11, 5).Value = 50
$ws.Cells.Item(11, 6).Value = 1.5
$ws.Cells.Item(11, 7).Value = 22.22
$ws.Cells.Item(11, 8).Value = 'Normal weight - Maintain a balanced diet and exercise.'
$ws.Cells.Item(11, 9).Value = 'Non-veg'
$ws.Cells.Item(11, 10).Value = 'vitamin_C'
$ws.Cells.Item(11, 11).Value = @'

 Fruit
  - Acerola, (west indian cherry), raw
  - Guavas, common, raw
  - Mango, Ataulfo, peeled, raw
  - Currants, european black, raw
  - Kiwifruit, ZESPRI SunGold, raw
  - Lemon peel, raw
  - Orange peel, raw
  - Longans, raw
  - Litchis, raw
  - Jujube, raw
  - Kiwifruit, green, raw
  - Persimmons, native, raw
  - Pummelo, raw
  - Strawberries, raw
  - Oranges, raw, with peel
  - Papayas, raw
  - Strawberries, raw
  - Kiwifruit (kiwi), green, peeled, raw
  - Pineapple, raw
 Juice
  - Acerola juice, raw
 Non Alcoholic
  - Beverages, tea, green, instant, decaffeinated, lemon, unsweetened, fortified with vitamin C
  - Beverages, tea, instant, lemon, with added ascorbic acid
 Vegetable
  - Peppers, hot chili, green, raw
  - Peppers, sweet, yellow, raw
  - Peppers, bell, orange, raw
  - Drumstick pods, raw
  - Pokeberry shoots, (poke), raw
  - Peppers, bell, yellow, raw
  - Peppers, bell, red, raw
  - Peppers, bell, green, raw
  - Balsam-pear (bitter gourd), pods, raw
  - Peppers, sweet, red, raw
  - Broccoli, flower clusters, raw
  - Broccoli, leaves, raw
  - Broccoli, stalks, raw
  - Peppers, sweet, green, raw
  - Taro, tahitian, raw
  - Broccoli, raw
  - Mustard spinach, (tendergreen), raw
  - Cauliflower, green, raw
'@
$ws.Cells.Item(11, 1).NumberFormat = $ws.Cells.Item(9, 1).NumberFormat

# --- New row 12 ---
$ws.Cells.Item(12, 1).Value = 45733.4061382423
$ws.Cells.Item(12, 2).Value = 'test'
$ws.Cells.Item(12, 3).Value = 25
$ws.Cells.Item(12, 4).Value = 'Male'
$ws.Cells.Item(12, 5).Value = 50
$ws.Cells.Item(12, 6).Value = 1.5
$ws.Cells.Item(12, 7).Value = 22.22
$ws.Cells.Item(12, 8).Value = 'Normal weight - Maintain a balanced diet and exercise.'
$ws.Cells.Item(12, 9).Value = 'Veg'
$ws.Cells.Item(12, 10).Value = 'vitamin_C'
$ws.Cells.Item(12, 11).Value = @'

 Fruit
  - Acerola, (west indian cherry), raw
  - Guavas, common, raw
  - Mango, Ataulfo, peeled, raw
  - Currants, european black, raw
  - Kiwifruit, ZESPRI SunGold, raw
  - Lemon peel, raw
  - Orange peel, raw
  - Longans, raw
  - Litchis, raw
  - Jujube, raw
  - Kiwifruit, green, raw
  - Persimmons, native, raw
  - Pummelo, raw
  - Strawberries, raw
  - Oranges, raw, with peel
  - Papayas, raw
  - Strawberries, raw
  - Kiwifruit (kiwi), green, peeled, raw
  - Pineapple, raw
 Juice
  - Acerola juice, raw
 Vegetable
  - Peppers, hot chili, green, raw
  - Peppers, sweet, yellow, raw
  - Peppers, bell, orange, raw
  - Drumstick pods, raw
  - Pokeberry shoots, (poke), raw
  - Peppers, bell, yellow, raw
  - Peppers, bell, red, raw
  - Peppers, bell, green, raw
  - Balsam-pear (bitter gourd), pods, raw
  - Peppers, sweet, red, raw
  - Broccoli, flower clusters, raw
  - Broccoli, leaves, raw
  - Broccoli, stalks, raw
  - Peppers, sweet, green, raw
  - Taro, tahitian, raw
  - Broccoli, raw
  - Mustard spinach, (tendergreen), raw
  - Cauliflower, green, raw
'@
$ws.Cells.Item(12, 1).NumberFormat = $ws.Cells.Item(9, 1).NumberFormat

